$d = $word.ActiveDocument

# Remove the floating text box ("Text Box 1") that currently shows
# "{statesAttorneyName}" above the first signature line.
$shp = $d.Shapes.Item(1)
$shp.Delete()

# Locate the paragraph that contains the two underscore signature
# lines ("___...___" TAB "___...___") that used to sit just below the
# text box, and insert a new paragraph right after it.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "_______________________________________" + [char]9 + "_______________________________________") {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIndex)
$p.Range.InsertParagraphAfter()

$newP = $d.Paragraphs.Item($targetIndex + 1)
$r = $newP.Range
$r.Text = "{statesAttorneyName}"

# Split the text into two runs ("{" and "statesAttorneyName}") to
# mirror the authored markup, without altering visible formatting.
$start = $newP.Range.Start
$r1 = $d.Range($start, $start + 1)
$r1.Font.Bold = 1
$r1.Font.Bold = 0
